$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.125.74'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -4.44%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.650.57'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.53%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.42'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.10%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5120'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.27%  '
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2589'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.68%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06428'
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07756'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.78%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.651.41'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.33%  '
$ws.Range("E13").Value = '  -4.78%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.879.02'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.46%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5514'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -5.81%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅7996'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.57%  '
$ws.Range("E17").Value = '  -5.71%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '26.155.70'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -4.39%  '
$ws.Range("E19").Value = '  +0.21%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '210.45'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.33%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.393'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.36%  '
$ws.Range("E22").Value = '  -4.02%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.041'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.41%  '
$ws.Range("E24").Value = '  +0.12%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.33'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.46%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.746'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.53%  '
$ws.Range("E27").Value = '  -2.88%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.966'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.68%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.78'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.87%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05131'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.79%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.242'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.85%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.349'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.46%  '
$ws.Range("E33").Value = '  -6.24%  '
$ws.Range("E34").Value = '  -4.60%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.738'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.77%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.358'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.52%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9223'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.27%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.169.09'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.95%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5688'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.79%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01585'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.22%  '
$ws.Range("E41").Value = '  +0.10%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.551'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.61%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.655'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.32%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8221'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.04%  '
$ws.Range("E45").Value = '  -0.93%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.789.33'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.45%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0₈117'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.68%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4548'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.43%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '55.44'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.68%  '
$ws.Range("E50").Value = '  +0.34%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.843'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.34%  '
